# Safety check - pre fetch from Dev.
#
# TestUsers.xlsx: add a "Manager" role for RickG (B2), mirror the
# "Normal User" role + original password onto MaggieG's row (B7/D7),
# and drop the stale hyperlink + its now-unused "Hyperlink" cell style
# that used to live on C2 (P@ssw0rd1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 carried an external mailto: hyperlink (and the built-in Hyperlink
# style that comes with it) - remove the link itself and restore plain
# "Normal" formatting to the cell, then drop the now-orphaned named
# style definition from the workbook.
$ws.Range("C2").Hyperlinks.Delete()
$ws.Range("C2").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

# New data: give RickG (row 2) a Manager role, and give MaggieG (row 7)
# the Normal User role with her original password repeated in the
# "New Pwd" column, matching row 8's pattern.
$ws.Range("B2").Value = "Manager"
$ws.Range("B7").Value = "Normal User"
$ws.Range("D7").Value = "P@ssw0rd1"

# Leave the selection on B2, the cell that was just edited.
$ws.Range("B2").Select()

$wb.Save()
